# Updated cryptos list with GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.848.10"
$ws.Range("E2").Value = "  -5.68%  "

# Row 3
$ws.Range("D3").Value = "1.819.34"
$ws.Range("E3").Value = "  -4.50%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.43%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.35"
$ws.Range("E5").Value = "  -2.64%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.36%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4621"
$ws.Range("E7").Value = "  -3.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3845"
$ws.Range("E8").Value = "  -3.95%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.84"
$ws.Range("E9").Value = "  -3.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07837"
$ws.Range("E10").Value = "  -2.53%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9589"
$ws.Range("E11").Value = "  -3.23%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.74"
$ws.Range("E12").Value = "  -6.33%  "

# Row 13
$ws.Range("D13").Value = "1.837.17"
$ws.Range("E13").Value = "  -3.69%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.629"
$ws.Range("E14").Value = "  -4.90%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.842"
$ws.Range("E15").Value = "  -3.81%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06857"
$ws.Range("E16").Value = "  +0.29%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.53%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.55"
$ws.Range("E18").Value = "  -2.90%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009917"
$ws.Range("E19").Value = "  -2.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.62"
$ws.Range("E20").Value = "  -4.35%  "

# Row 21
$ws.Range("E21").Value = "  -0.40%  "

# Row 22
$ws.Range("D22").Value = "27.882.41"
$ws.Range("E22").Value = "  -5.60%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.299"
$ws.Range("E23").Value = "  -3.76%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  -5.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.098"
$ws.Range("E25").Value = "  -2.88%  "

# Row 26
$ws.Range("D26").Value = "2.027.14"
$ws.Range("E26").Value = "  -5.22%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.19"
$ws.Range("E27").Value = "  -3.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.16"
$ws.Range("E28").Value = "  -1.90%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.709"
$ws.Range("E29").Value = "  -11.87%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.963"
$ws.Range("E30").Value = "  -4.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.32"
$ws.Range("E31").Value = "  -2.25%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09235"
$ws.Range("E32").Value = "  -3.08%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9328"
$ws.Range("E33").Value = "  -6.41%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.270"
$ws.Range("E34").Value = "  -3.62%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.405"
$ws.Range("E35").Value = "  -3.82%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.311"
$ws.Range("E36").Value = "  -5.50%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05924"
$ws.Range("E37").Value = "  -8.30%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02138"
$ws.Range("E38").Value = "  -4.46%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.144"
$ws.Range("E39").Value = "  -3.93%  "

# Row 40
$ws.Range("E40").Value = "  -0.43%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.567"
$ws.Range("E41").Value = "  -2.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5554"
$ws.Range("E42").Value = "  -4.49%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.879"
$ws.Range("E43").Value = "  -6.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1763"
$ws.Range("E44").Value = "  -3.15%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.244"
$ws.Range("E45").Value = "  +0.65%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.219"
$ws.Range("E46").Value = "  -9.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.61"
$ws.Range("E47").Value = "  -4.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5229"
$ws.Range("E48").Value = "  -4.55%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06982"
$ws.Range("E49").Value = "  -5.91%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.815"
$ws.Range("E50").Value = "  -6.76%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.28"
$ws.Range("E51").Value = "  -3.23%  "
